$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "g/d"
$ws.Range("C3").Value = "g/d"
$ws.Range("C4").Value = "g/d"
$ws.Range("C5").Value = "g/d"
$ws.Range("C6").Value = "Kcal/d"
$ws.Range("C7").Value = "g/d"
$ws.Range("C8").Value = "g/d"
$ws.Range("C9").Value = "mg/d"
$ws.Range("C10").Value = "mg/d"
$ws.Range("C11").Value = "mg/d"
$ws.Range("C12").Value = "IU/d"
$ws.Range("C13").Value = "mg/d retinol"
$ws.Range("C14").Value = "ug/d"

$ws.Range("D8").Select() | Out-Null
